# GameShop.xlsx update: "add some item to recover health and mental. quest of 1-2"
#  - Row 40 (ItemId 22034011 / "体力药水") becomes ItemId 22034010 / "饼干"
#  - Two brand-new rows are appended to the 表3 table:
#       row 41: Id 15000040, ItemId 22034011, Shelf 3, ~Name "红色胶囊"
#       row 42: Id 15000041, ItemId 22034012, Shelf 3, ~Name "蓝色胶囊"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

$formula = "=LOOKUP(表3[[#This Row],[ItemId]],[1]其他!`$A:`$A,[1]其他!`$B:`$B)"

# --- update the existing last row (row 40): ItemId changes, ~Name formula unchanged ---
$ws.Range("B40").Value = 22034010

# --- append row 41 ---
$row41 = $lo.ListRows.Add()
$row41.Range.Cells.Item(1,1).Value = 15000040
$row41.Range.Cells.Item(1,2).Value = 22034011
$row41.Range.Cells.Item(1,3).Value = 3
$row41.Range.Cells.Item(1,4).Formula = $formula

# --- append row 42 ---
$row42 = $lo.ListRows.Add()
$row42.Range.Cells.Item(1,1).Value = 15000041
$row42.Range.Cells.Item(1,2).Value = 22034012
$row42.Range.Cells.Item(1,3).Value = 3
$row42.Range.Cells.Item(1,4).Formula = $formula

# --- formatting: calculated column picks up D40's style (thin border, no fill) ---
$ws.Range("D40").Copy()
$row41.Range.Cells.Item(1,4).PasteSpecial(-4122)
$ws.Range("D40").Copy()
$row42.Range.Cells.Item(1,4).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- sheet view matches the authored selection after the edit ---
$ws.Range("D31").Select()
